# ---------------------------------------------------------------------
# Rename sheet 1, and create sheet 2 as a literal copy of sheet 1 so it
# starts out with identical sheetPr / pageMargins / sheetFormatPr, then
# wipe sheet 2's cell contents (keeping its page/sheet metadata intact).
# ---------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Clear()

$ws1.Name = "Data Parkir"
$ws2.Name = "Ringkasan Keuangan"

# ---------------------------------------------------------------------
# Sheet 1: "Data Parkir"
# ---------------------------------------------------------------------
# A1 already carries the bold/bordered/centered header style from the
# original template (style index 1) - propagate its *format only* across
# the rest of the header row before we touch any values.
$ws1.Range("A1").Copy()
$ws1.Range("B1:K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headers1 = @("Kode_Parking","No_Kendaraan","Jenis_Kendaraan","Waktu_Masuk","Waktu_Keluar","Durasi","Biaya","Uang_Pembayaran","Nama_Petugas","Foto_Masuk","Foto_Keluar")
for ($c = 0; $c -lt $headers1.Length; $c++) {
    $ws1.Cells.Item(1, $c + 1).Value = $headers1[$c]
}

$rows1 = @(
    @("NJ9E1D0UMB22","Y 0313 HAS","Motor","2025-01-31 21:09:22","2025-01-31 21:09:47","00:00:25",2000,5000,"Reza Ramdan Permana","./capture/masuk/NJ9E1D0UMB22.png","./capture/keluar/NJ9E1D0UMB22.png"),
    @("SN12XKBEG18L","D 4230 ASQ","Mobil","2025-01-31 21:10:39","2025-01-31 21:11:58","00:01:19",4000,60000,"Reza Ramdan Permana","./capture/masuk/SN12XKBEG18L.png","./capture/keluar/SN12XKBEG18L.png"),
    @("58O0J9BUUNER","D 9530 JFD","Motor","2025-01-31 21:15:34","2025-01-31 21:16:09","00:00:35",2000,5000,"Reza Ramdan Permana","./capture/masuk/58O0J9BUUNER.png","./capture/keluar/58O0J9BUUNER.png")
)

for ($r = 0; $r -lt $rows1.Length; $r++) {
    $row = $rows1[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------
# Sheet 2: "Ringkasan Keuangan"
# ---------------------------------------------------------------------
# $ws2.Cells.Clear() above also wiped A1's inherited header format, so
# re-stamp it from the still-formatted $ws1.Range("A1") before fanning
# the format out across the rest of the header row.
$ws1.Range("A1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws2.Cells.Item(1, 1).Value = "Deskripsi"
$ws2.Cells.Item(1, 2).Value = "Nilai"

$rows2 = @(
    @("Total Pemasukan", "Rp 8,000"),
    @("Total Kembalian", "Rp 62,000"),
    @("Uang yang Perlu Disetorkan", "Rp 8,000"),
    @("Terbilang", "Delapan ribu Rupiah")
)

for ($r = 0; $r -lt $rows2.Length; $r++) {
    $row = $rows2[$r]
    $ws2.Cells.Item($r + 2, 1).Value = $row[0]
    $ws2.Cells.Item($r + 2, 2).Value = $row[1]
}

# Restore selection/active state to match the original layout.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
